$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row => [D value, E value] (only cells that actually change per the diff).
# $null means "this column is unchanged for this row".
$updates = @{
    2  = @("68.462.15", "  +1.12%  ")
    3  = @("3.748.01",  "  -0.87%  ")
    4  = @($null,       "  +0.10%  ")
    5  = @("595.71",    "  -0.18%  ")
    6  = @("167.49",    "  -0.81%  ")
    7  = @("3.746.03",  "  -0.83%  ")
    8  = @($null,       "  -0.10%  ")
    10 = @($null,       "  -2.73%  ")
    11 = @($null,       "  -0.35%  ")
    12 = @($null,       "  -1.22%  ")
    13 = @($null,       "  -6.34%  ")
    14 = @("36.06",     "  -0.85%  ")
    15 = @("4.378.24",  "  -0.86%  ")
    16 = @("3.754.13",  "  -1.11%  ")
    17 = @("68.456.64", "  +1.17%  ")
    18 = @("17.96",     "  -3.23%  ")
    19 = @($null,       "  -2.49%  ")
    20 = @($null,       "  -0.12%  ")
    21 = @($null,       "  +2.36%  ")
    22 = @("465.36",    "  -0.58%  ")
    23 = @("0.698",     "  -2.71%  ")
    24 = @("84.58",     "  +0.93%  ")
    25 = @($null,       "  -1.88%  ")
    26 = @($null,       "  -0.65%  ")
    27 = @("11.99",     "  -1.15%  ")
    28 = @($null,       "  -0.07%  ")
    29 = @($null,       "  -3.18%  ")
    30 = @("3.894.36",  "  -0.97%  ")
    31 = @($null,       "  -4.43%  ")
    32 = @("7.31",      "  -4.03%  ")
    33 = @("29.82",     "  -2.42%  ")
    34 = @($null,       "  -1.99%  ")
    35 = @($null,       "  +0.91%  ")
    37 = @("3.702.29",  "  -1.17%  ")
    38 = @($null,       "  -2.65%  ")
    39 = @("3.36",      "  -9.72%  ")
    40 = @($null,       "  +0.89%  ")
    41 = @("0.999",     "  -0.30%  ")
    42 = @("5.82",      "  +0.16%  ")
    43 = @($null,       "  +0.06%  ")
    45 = @($null,       "  -2.26%  ")
    46 = @("43.81",     "  +11.34%  ")
    47 = @($null,       "  -0.99%  ")
    48 = @($null,       "  -0.80%  ")
    49 = @("45.99",     "  +0.58%  ")
    50 = @("146.77",    "  +4.39%  ")
    51 = @("389.93",    "  -1.51%  ")
}

# Column D cells whose new text would otherwise be auto-recognised as a plain
# number by Excel (e.g. "595.71"). The source file stores every Price/Volume
# cell as text, so for these rows we force the Text number format before
# writing, then restore the "Normal" style afterwards so no stray style
# index is left behind on the cell.
$forceTextRows = @(5, 6, 14, 18, 22, 23, 24, 27, 32, 33, 39, 41, 42, 46, 49, 50, 51)

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]

    if ($dVal -ne $null) {
        $dCell = $ws.Cells.Item($row, 4)
        if ($forceTextRows -contains $row) {
            $dCell.NumberFormat = "@"
            $dCell.Value = $dVal
            $dCell.Style = "Normal"
        } else {
            $dCell.Value = $dVal
        }
    }

    $ws.Cells.Item($row, 5).Value = $eVal
}
